$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# Delete entire columns I:K (INDICE, CLUSTER, RANKING)
$ws.Range("I1:K1").EntireColumn.Delete()
